$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 8 new rows above the current row 2 (shifts existing data rows 2-21 down to 10-29)
$ws.Range("A2:C9").Insert(-4121)
# The inserted cells inherit formatting from the row above (the bold header row) -
# strip that back off so the new data rows are plain, like the rest of the data rows.
$ws.Range("A2:C9").ClearFormats()

# Fill in the new rows 2-9 with the newly-added data
$newRows = @(
  @(-0.007177666760981, 0.0068722339347004, 0.0154243474826216),
  @(-0.0125227374956011, 0.0532979927957057, -0.0287106670439243),
  @(-0.0233655963093042, 0.0858265683054924, 0.0064140851609408),
  @(-0.0406225398182868, -0.0181732401251792, 0.0048869219608604),
  @(0.0109955742955207, -0.015118914656341, 0.0439822971820831),
  @(0.0197004042565822, -0.0032070425804704, 0.0224492978304624),
  @(0.0274889357388019, 0.0001527163112768, 0.0108428578823804),
  @(0.0339030213654041, 0.0074830991216003, -0.0473420582711696)
)

$r = 2
foreach ($row in $newRows) {
  $ws.Cells.Item($r, 1).Value = $row[0]
  $ws.Cells.Item($r, 2).Value = $row[1]
  $ws.Cells.Item($r, 3).Value = $row[2]
  $r = $r + 1
}

# Append 2 new rows at the bottom (rows 30-31)
$appendRows = @(
  @(-0.0207694191485643, -0.0707076489925384, -0.09605856239795681),
  @(0.0001527163112768, -0.0478002056479454, -0.06643159687519069)
)

$r = 30
foreach ($row in $appendRows) {
  $ws.Cells.Item($r, 1).Value = $row[0]
  $ws.Cells.Item($r, 2).Value = $row[1]
  $ws.Cells.Item($r, 3).Value = $row[2]
  $r = $r + 1
}
